$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates (puzzle-area rework in lvl3 map)
$ws.Cells.Item(1, 36).Value = 11
$ws.Cells.Item(2, 36).Value = 11
$ws.Cells.Item(3, 36).Value = 11
$ws.Cells.Item(4, 36).Value = 11
$ws.Cells.Item(5, 36).Value = 11
$ws.Cells.Item(10, 29).Value = 5
$ws.Cells.Item(14, 34).Value = 17
$ws.Cells.Item(14, 38).Value = 17
$ws.Cells.Item(15, 34).Value = 17
$ws.Cells.Item(15, 38).Value = 17
$ws.Cells.Item(16, 34).Value = 17
$ws.Cells.Item(16, 38).Value = 17
$ws.Cells.Item(17, 34).Value = 17
$ws.Cells.Item(17, 38).Value = 17
$ws.Cells.Item(24, 8).Value = 16
$ws.Cells.Item(36, 25).Value = 16
$ws.Cells.Item(37, 26).Value = 5
$ws.Cells.Item(37, 27).Value = 5
$ws.Cells.Item(41, 36).Value = 0
$ws.Cells.Item(44, 37).Value = 7
$ws.Cells.Item(45, 36).Value = 10
$ws.Cells.Item(46, 36).Value = 7
$ws.Cells.Item(47, 37).Value = 10
$ws.Cells.Item(49, 36).Value = 10
$ws.Cells.Item(50, 37).Value = 0
$ws.Cells.Item(51, 2).Value = 12
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = 12
$ws.Cells.Item(52, 37).Value = 0
$ws.Cells.Item(52, 39).Value = 7
$ws.Cells.Item(53, 1).Value = 2
$ws.Cells.Item(53, 2).Value = 12
$ws.Cells.Item(53, 12).Value = 7
$ws.Cells.Item(53, 13).Value = 7
$ws.Cells.Item(53, 14).Value = 7
$ws.Cells.Item(53, 15).Value = 7
$ws.Cells.Item(53, 16).Value = 7
$ws.Cells.Item(53, 17).Value = 7
$ws.Cells.Item(53, 18).Value = 7
$ws.Cells.Item(53, 19).Value = 7
$ws.Cells.Item(53, 20).Value = 7
$ws.Cells.Item(53, 29).Value = 16
$ws.Cells.Item(54, 2).Value = 12
$ws.Cells.Item(54, 12).Value = 10
$ws.Cells.Item(54, 13).Value = 8
$ws.Cells.Item(54, 19).Value = 8
$ws.Cells.Item(54, 20).Value = 10
$ws.Cells.Item(55, 2).Value = 12
$ws.Cells.Item(55, 12).Value = 10
$ws.Cells.Item(55, 13).Value = 8
$ws.Cells.Item(55, 16).Value = 14
$ws.Cells.Item(55, 19).Value = 8
$ws.Cells.Item(55, 20).Value = 10
$ws.Cells.Item(56, 2).Value = 12
$ws.Cells.Item(56, 12).Value = 10
$ws.Cells.Item(56, 13).Value = 7
$ws.Cells.Item(56, 14).Value = 7
$ws.Cells.Item(56, 15).Value = 7
$ws.Cells.Item(56, 16).Value = 7
$ws.Cells.Item(56, 17).Value = 7
$ws.Cells.Item(56, 18).Value = 7
$ws.Cells.Item(56, 19).Value = 7
$ws.Cells.Item(56, 20).Value = 10
$ws.Cells.Item(57, 12).Value = 10
$ws.Cells.Item(57, 20).Value = 10
$ws.Cells.Item(58, 12).Value = 7
$ws.Cells.Item(58, 14).Value = 18
$ws.Cells.Item(58, 18).Value = 18
$ws.Cells.Item(58, 20).Value = 7
$ws.Cells.Item(59, 12).Value = 7
$ws.Cells.Item(59, 13).Value = 7
$ws.Cells.Item(59, 14).Value = 7
$ws.Cells.Item(59, 15).Value = 9
$ws.Cells.Item(59, 16).Value = 9
$ws.Cells.Item(59, 17).Value = 9
$ws.Cells.Item(59, 18).Value = 7
$ws.Cells.Item(59, 19).Value = 7
$ws.Cells.Item(59, 20).Value = 7
$ws.Cells.Item(60, 13).Value = 8
$ws.Cells.Item(60, 19).Value = 8
$ws.Cells.Item(61, 11).Value = 20
$ws.Cells.Item(61, 13).Value = 8
$ws.Cells.Item(61, 15).Value = 19
$ws.Cells.Item(61, 17).Value = 19
$ws.Cells.Item(61, 19).Value = 8
$ws.Cells.Item(61, 21).Value = 20
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 1
$ws.Cells.Item(62, 11).Value = 1
$ws.Cells.Item(62, 12).Value = 9
$ws.Cells.Item(62, 13).Value = 7
$ws.Cells.Item(62, 14).Value = 7
$ws.Cells.Item(62, 15).Value = 7
$ws.Cells.Item(62, 16).Value = 7
$ws.Cells.Item(62, 17).Value = 7
$ws.Cells.Item(62, 18).Value = 9
$ws.Cells.Item(62, 19).Value = 7
$ws.Cells.Item(62, 20).Value = 9
$ws.Cells.Item(62, 21).Value = 1
$ws.Cells.Item(62, 22).Value = 1
$ws.Cells.Item(62, 23).Value = 1
$ws.Cells.Item(62, 24).Value = 1
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = 2
$ws.Cells.Item(63, 11).Value = 2
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 18
$ws.Cells.Item(63, 14).Value = 0
$ws.Cells.Item(63, 15).Value = 0
$ws.Cells.Item(63, 16).Value = 0
$ws.Cells.Item(63, 17).Value = 0
$ws.Cells.Item(63, 18).Value = 0
$ws.Cells.Item(63, 19).Value = 0
$ws.Cells.Item(63, 20).Value = 0
$ws.Cells.Item(63, 21).Value = 2
$ws.Cells.Item(63, 22).Value = 2
$ws.Cells.Item(63, 23).Value = 2
$ws.Cells.Item(63, 24).Value = 2
$ws.Cells.Item(64, 12).Value = 7
$ws.Cells.Item(64, 13).Value = 7
$ws.Cells.Item(64, 17).Value = 7
$ws.Cells.Item(64, 18).Value = 7
$ws.Cells.Item(64, 19).Value = 7
$ws.Cells.Item(64, 20).Value = 7
$ws.Cells.Item(66, 36).Value = 3
$ws.Cells.Item(67, 36).Value = 3
$ws.Cells.Item(68, 36).Value = 3
$ws.Cells.Item(69, 36).Value = 3
$ws.Cells.Item(70, 36).Value = 3
$ws.Cells.Item(75, 29).Value = 1
$ws.Cells.Item(79, 34).Value = 4
$ws.Cells.Item(79, 38).Value = 4
$ws.Cells.Item(80, 34).Value = 4
$ws.Cells.Item(80, 38).Value = 4
$ws.Cells.Item(81, 34).Value = 4
$ws.Cells.Item(81, 38).Value = 4
$ws.Cells.Item(82, 34).Value = 4
$ws.Cells.Item(82, 38).Value = 4
$ws.Cells.Item(102, 26).Value = 1
$ws.Cells.Item(102, 27).Value = 1
$ws.Cells.Item(106, 36).Value = 0
$ws.Cells.Item(109, 37).Value = 1
$ws.Cells.Item(110, 36).Value = 3
$ws.Cells.Item(111, 36).Value = 1
$ws.Cells.Item(112, 37).Value = 3
$ws.Cells.Item(114, 36).Value = 3
$ws.Cells.Item(115, 37).Value = 0
$ws.Cells.Item(116, 2).Value = 3
$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = 3
$ws.Cells.Item(117, 37).Value = 0
$ws.Cells.Item(117, 39).Value = 1
$ws.Cells.Item(118, 1).Value = 1
$ws.Cells.Item(118, 2).Value = 3
$ws.Cells.Item(118, 12).Value = 1
$ws.Cells.Item(118, 13).Value = 1
$ws.Cells.Item(118, 14).Value = 1
$ws.Cells.Item(118, 15).Value = 1
$ws.Cells.Item(118, 16).Value = 1
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = 1
$ws.Cells.Item(118, 19).Value = 1
$ws.Cells.Item(118, 20).Value = 1
$ws.Cells.Item(119, 2).Value = 3
$ws.Cells.Item(119, 12).Value = 3
$ws.Cells.Item(119, 13).Value = 1
$ws.Cells.Item(119, 19).Value = 1
$ws.Cells.Item(119, 20).Value = 3
$ws.Cells.Item(120, 2).Value = 3
$ws.Cells.Item(120, 12).Value = 3
$ws.Cells.Item(120, 13).Value = 1
$ws.Cells.Item(120, 19).Value = 1
$ws.Cells.Item(120, 20).Value = 3
$ws.Cells.Item(121, 2).Value = 3
$ws.Cells.Item(121, 12).Value = 3
$ws.Cells.Item(121, 13).Value = 1
$ws.Cells.Item(121, 14).Value = 1
$ws.Cells.Item(121, 15).Value = 1
$ws.Cells.Item(121, 16).Value = 1
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = 1
$ws.Cells.Item(121, 19).Value = 1
$ws.Cells.Item(121, 20).Value = 3
$ws.Cells.Item(122, 12).Value = 3
$ws.Cells.Item(122, 20).Value = 3
$ws.Cells.Item(123, 12).Value = 1
$ws.Cells.Item(123, 14).Value = 1
$ws.Cells.Item(123, 18).Value = 1
$ws.Cells.Item(123, 20).Value = 1
$ws.Cells.Item(124, 12).Value = 1
$ws.Cells.Item(124, 13).Value = 1
$ws.Cells.Item(124, 14).Value = 1
$ws.Cells.Item(124, 15).Value = 2
$ws.Cells.Item(124, 16).Value = 2
$ws.Cells.Item(124, 17).Value = 2
$ws.Cells.Item(124, 18).Value = 1
$ws.Cells.Item(124, 19).Value = 1
$ws.Cells.Item(124, 20).Value = 1
$ws.Cells.Item(125, 13).Value = 1
$ws.Cells.Item(125, 19).Value = 1
$ws.Cells.Item(126, 13).Value = 1
$ws.Cells.Item(126, 19).Value = 1
$ws.Cells.Item(127, 9).Value = 1
$ws.Cells.Item(127, 10).Value = 1
$ws.Cells.Item(127, 11).Value = 1
$ws.Cells.Item(127, 12).Value = 2
$ws.Cells.Item(127, 13).Value = 1
$ws.Cells.Item(127, 14).Value = 1
$ws.Cells.Item(127, 15).Value = 1
$ws.Cells.Item(127, 16).Value = 1
$ws.Cells.Item(127, 17).Value = 1
$ws.Cells.Item(127, 18).Value = 2
$ws.Cells.Item(127, 19).Value = 1
$ws.Cells.Item(127, 20).Value = 2
$ws.Cells.Item(127, 21).Value = 1
$ws.Cells.Item(127, 22).Value = 1
$ws.Cells.Item(127, 23).Value = 1
$ws.Cells.Item(127, 24).Value = 1
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 14).Value = 0
$ws.Cells.Item(128, 15).Value = 0
$ws.Cells.Item(128, 16).Value = 0
$ws.Cells.Item(128, 17).Value = 0
$ws.Cells.Item(128, 18).Value = 0
$ws.Cells.Item(128, 19).Value = 0
$ws.Cells.Item(128, 20).Value = 0

# Selection / view state
$ws.Range("L63").Select()
